$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.741684079170227
$ws.Range("B1").Value = 1.880554556846619
$ws.Range("C1").Value = 1.907991170883179
$ws.Range("D1").Value = 2.484105587005615
$ws.Range("E1").Value = 2.834566593170166
